$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (every cell in the workbook that shows this status gets the new text)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 2. Latest Handback DateTime (column H) now has a real timestamp instead of
#    the "0001-01-01 00:00:00" placeholder, per locale sheet.
# ---------------------------------------------------------------------------
$wsZhCn.Range("H2").Value = "2016-03-18 07:28:12"
$wsZhCn.Range("H3").Value = "2016-03-18 07:28:12"

$wsDeDe.Range("H2").Value = "2016-03-18 07:28:16"
$wsDeDe.Range("H3").Value = "2016-03-18 07:28:16"

# ---------------------------------------------------------------------------
# 3. Populate "Latest Target File" (F) and "Latest Handback File" (G) columns
#    with hyperlinked file names, mirroring the existing Source File Name (A)
#    and Latest Handoff File (D) links for the same row.
# ---------------------------------------------------------------------------

# zh-cn, row 2 (416ce689-...)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/972fa0d7308bcac74b74e1b12093337f5156661e/e2e/416ce689-5ba7-4604-92bd-dd0b924fa3fe.md", "", "", "416ce689-5ba7-4604-92bd-dd0b924fa3fe.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5e91501aa781806f5db99307d0338807048376c3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/416ce689-5ba7-4604-92bd-dd0b924fa3fe.127335717d341eb863826d100b06fd677421f8b1.zh-cn.xlf", "", "", "416ce689-5ba7-4604-92bd-dd0b924fa3fe.127335717d341eb863826d100b06fd677421f8b1.zh-cn.xlf")

# zh-cn, row 3 (4ffe8906-...)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/972fa0d7308bcac74b74e1b12093337f5156661e/e2e/4ffe8906-cf1d-4fba-866d-8e285fb1ae43.md", "", "", "4ffe8906-cf1d-4fba-866d-8e285fb1ae43.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5e91501aa781806f5db99307d0338807048376c3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4ffe8906-cf1d-4fba-866d-8e285fb1ae43.185f710f20f2bc8c2e415ddc53f41093cbddac02.zh-cn.xlf", "", "", "4ffe8906-cf1d-4fba-866d-8e285fb1ae43.185f710f20f2bc8c2e415ddc53f41093cbddac02.zh-cn.xlf")

# de-de, row 2 (416ce689-...)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/972fa0d7308bcac74b74e1b12093337f5156661e/e2e/416ce689-5ba7-4604-92bd-dd0b924fa3fe.md", "", "", "416ce689-5ba7-4604-92bd-dd0b924fa3fe.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f6a015c2bdcd279394a3f94551ab7819ed067a88/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/416ce689-5ba7-4604-92bd-dd0b924fa3fe.127335717d341eb863826d100b06fd677421f8b1.de-de.xlf", "", "", "416ce689-5ba7-4604-92bd-dd0b924fa3fe.127335717d341eb863826d100b06fd677421f8b1.de-de.xlf")

# de-de, row 3 (4ffe8906-...)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/972fa0d7308bcac74b74e1b12093337f5156661e/e2e/4ffe8906-cf1d-4fba-866d-8e285fb1ae43.md", "", "", "4ffe8906-cf1d-4fba-866d-8e285fb1ae43.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f6a015c2bdcd279394a3f94551ab7819ed067a88/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4ffe8906-cf1d-4fba-866d-8e285fb1ae43.185f710f20f2bc8c2e415ddc53f41093cbddac02.de-de.xlf", "", "", "4ffe8906-cf1d-4fba-866d-8e285fb1ae43.185f710f20f2bc8c2e415ddc53f41093cbddac02.de-de.xlf")

# Match the look of the existing hyperlink cells (underlined, blue) on the
# newly-added cells.
$newLinkCells = @("F2", "G2", "F3", "G3")
foreach ($addr in $newLinkCells) {
    $wsZhCn.Range($addr).Font.Underline = $true
    $wsZhCn.Range($addr).Font.Color = 0xED9564
    $wsDeDe.Range($addr).Font.Underline = $true
    $wsDeDe.Range($addr).Font.Color = 0xED9564
}
